$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.35, 701999.4142380001, 25485.65933, 676513.754908, 15696.92096, 327052.414802, 16903.008564, 310149.406238, 17252.32464266667),
    @(0.4, 2301811.046997, 31830.510612, 2269980.536385, 19497.82999666667, 987316.1288909999, 19039.975669, 968276.153222, 21175.29658233333),
    @(0.45, 2574588.917021, 36607.053794, 2537981.863227, 23046.86977733333, 1025262.224689, 20152.894303, 1005109.330386, 24845.11606266667)
)

$startRow = 10
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
